$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "272.82"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.06%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.79"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.36%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.907"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.09%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06323"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.91%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.904"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.36%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.361"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "5.55%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.326"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "45.11%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8871"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3.08%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.10%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.05152"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.79%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07388"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.24%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03133"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.49%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09044"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.02%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001559"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.78%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006319"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3.86%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006025"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.51%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.464"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.22%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.31%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.62%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1335"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.19%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.913"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.00%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04355"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.59%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001177"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.43%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.003670"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-12.50%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.24%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001701"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "1.62%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04025"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.29%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006613"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "6.68%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.99%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002104"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.95%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.19%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005319"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.67%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.358"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "162.65%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-12.75%"
